# Daily update at 8 AM UTC
# Appends the next day's row (Day, Chase, Bryce, Zach) to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: next day after the existing 2025-10-21 (serial 45951) is
# 2025-10-22 (serial 45952), with that day's win counts for each player.
$ws.Range("A3").Value = 45952
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 4

# Re-apply the date/time formatting used by this tracker. The existing date
# cell picks up an explicit date-time format, and the freshly appended date
# cell gets the plain date format - matching the sheet's day-over-day style
# churn.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A3").NumberFormat = "yyyy-mm-dd"
$ws.Range("A3").NumberFormat = "YYYY-MM-DD"
